$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the COST value (currently in F2/F3) before it gets overwritten -
# it needs to move to the new I2/I3 "COST" column.
$cost2 = $ws.Range("F2").Value2
$cost3 = $ws.Range("F3").Value2

# Headers: F1 becomes "GENDER"; new columns G1/H1/I1 are added for
# CATEGORY / SUB CATEGORY / COST.
$ws.Range("F1").Value2 = "GENDER"
$ws.Range("G1").Value2 = "CATEGORY"
$ws.Range("H1").Value2 = "SUB CATEGORY"
$ws.Range("I1").Value2 = "COST"

# Row 2 data
$ws.Range("F2").Value2 = "M"
$ws.Range("G2").Value2 = "cat1"
$ws.Range("H2").Value2 = "sub1"
$ws.Range("I2").Value2 = $cost2

# Row 3 data
$ws.Range("F3").Value2 = "M"
$ws.Range("G3").Value2 = "cat2"
$ws.Range("H3").Value2 = "sub2"
$ws.Range("I3").Value2 = $cost3

# Columns F:I switch from the "General" numeric style to the "@" text
# style already used by column D (this reuses the existing style, matching
# the template/dropdown formatting applied to the new columns).
$ws.Range("F1:I3").NumberFormat = "@"

# Match the target workbook's active cell selection.
$ws.Range("H6").Select() | Out-Null
